# Fruta / hortaliza, semanal
# Insert a new weekly price-report row for "Ají" (Terminal Hortofrutícola
# Agro Chillán) above the current row 38, shifting the existing rows
# 38-47 down to 39-48.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 38:47 down one row to make room for the new record.
$ws.Rows.Item(38).Insert()

# Populate the newly inserted row 38 with the new weekly observation.
$ws.Range("A38").Value = 7
$ws.Range("B38").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C38").Value = 'Ñuble'
$ws.Range("D38").Value = 44551
$ws.Range("E38").Value = 16
$ws.Range("F38").Value = 100112021
$ws.Range("G38").Value = 'Ají'
$ws.Range("H38").Value = 'Americana (o)'
$ws.Range("I38").Value = 'Primera'
$ws.Range("J38").Value = 60
$ws.Range("K38").Value = 17500
$ws.Range("L38").Value = 18000
$ws.Range("M38").Value = 17750
$ws.Range("N38").Value = '$/caja 15 kilos'
$ws.Range("O38").Value = 'Región del Maule'
$ws.Range("P38").Value = 1183
$ws.Range("Q38").Value = 15
$ws.Range("R38").Value = 'Hortaliza'
